{"js": "// The \"Bonus:\" section contains a bulleted list of ideas. The edit:\n//  1. Turns the \"Create a 3+ frame explosion...\" bullet into \"Add the concept of lives\"\n//  2. Inserts a new bullet right after it with the old \"Alter the code...\" text\n//  3. Turns the (old) \"Add the concept of lives\" bullet into the \"start screen\" bullet\n//  4. Removes the now-duplicated old \"Alter the code...\" bullet\n//  5. Removes the now-duplicated old \"start screen\" bullet\n//  6. Turns the \"Add a timer...\" bullet into the new \"free play\" mode bullet\n//\n// insertOoxml(..., Word.InsertLocation.replace) gives full control over the\n// exact run/formatting structure, matching the target XML precisely. Each\n// Paragraph object handle becomes unusable for a further insertOoxml call\n// once it has already been used for one, so we always re-locate paragraphs\n// by their (currently still unique) text right before each such call.\n\nconst body = context.document.body;\n\nfunction findIndex(items, matchText) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text.indexOf(matchText) !== -1) {\n      return i;\n    }\n  }\n  throw new Error(\"Paragraph not found: \" + matchText);\n}\n\nfunction findAllIndices(items, matchText) {\n  const out = [];\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text.indexOf(matchText) !== -1) {\n      out.push(i);\n    }\n  }\n  return out;\n}\n\nfunction wrapPkg(pBodyXml) {\n  return (\n    '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    pBodyXml +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\nconst livesOoxml = wrapPkg(\n  '<w:p>' +\n  '<w:pPr><w:ind w:firstLine=\"540\"/><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr></w:pPr>' +\n  '<w:r><w:t xml:space=\"preserve\">\\u25a1 </w:t></w:r>' +\n  '<w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\">Add </w:t></w:r>' +\n  '<w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\">the concept of </w:t></w:r>' +\n  '<w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>lives</w:t></w:r>' +\n  '</w:p>'\n);\n\nconst shiesOoxml = wrapPkg(\n  '<w:p>' +\n  '<w:pPr><w:ind w:firstLine=\"540\"/><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr></w:pPr>' +\n  '<w:r><w:t>\\u25a1</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>Alter the code so instead of following a line, the tracking bot \\u201cshies\\u201d away from a line</w:t></w:r>' +\n  '</w:p>'\n);\n\nconst startScreenOoxml = wrapPkg(\n  '<w:p>' +\n  '<w:pPr><w:ind w:firstLine=\"540\"/><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\">        </w:t></w:r>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  '<w:r><w:t xml:space=\"preserve\">\\u25a1 </w:t></w:r>' +\n  '<w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> Create</w:t></w:r>' +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  '<w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> a start screen that allows the user to select which mode (follow or shy) to use</w:t></w:r>' +\n  '</w:p>'\n);\n\nconst freePlayOoxml = wrapPkg(\n  '<w:p>' +\n  '<w:pPr><w:ind w:firstLine=\"540\"/><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr></w:pPr>' +\n  '<w:r><w:t xml:space=\"preserve\">\\u25a1 </w:t></w:r>' +\n  '<w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>Create a \\u201cfree play\\u201d mode where the player gets points for hitting random popup sprites</w:t></w:r>' +\n  '</w:p>'\n);\n\n// --- Step 1: \"Create a 3+ frame explosion...\" -> \"Add the concept of lives\" ---\nbody.paragraphs.load(\"items,text\");\nawait context.sync();\nlet idx = findIndex(body.paragraphs.items, \"explosion animation\");\nbody.paragraphs.items[idx].insertOoxml(livesOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Step 2: insert a placeholder paragraph after it, then give it the \"Alter the code...\" content ---\nbody.paragraphs.load(\"items,text\");\nawait context.sync();\nidx = findIndex(body.paragraphs.items, \"the concept of lives\");\nbody.paragraphs.items[idx].insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\nbody.paragraphs.load(\"items,text\");\nawait context.sync();\nidx = findIndex(body.paragraphs.items, \"the concept of lives\");\nbody.paragraphs.items[idx + 1].insertOoxml(shiesOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Step 3: turn the ORIGINAL \"Add the concept of lives\" bullet (further down) into the start-screen bullet ---\nbody.paragraphs.load(\"items,text\");\nawait context.sync();\nlet livesIdxs = findAllIndices(body.paragraphs.items, \"the concept of lives\");\n// livesIdxs[0] is the bullet we created in Step 1 (keep it); the stale\n// original is the later duplicate.\nbody.paragraphs.items[livesIdxs[1]].insertOoxml(startScreenOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Step 4: delete the now-duplicated old \"Alter the code...\" bullet ---\nbody.paragraphs.load(\"items,text\");\nawait context.sync();\nlet shiesIdxs = findAllIndices(body.paragraphs.items, \"\\u201cshies\\u201d away from a line\");\n// shiesIdxs[0] is the bullet we created in Step 2 (keep it); delete the\n// later, now-stale duplicate.\nbody.paragraphs.items[shiesIdxs[1]].delete();\nawait context.sync();\n\n// --- Step 5: delete the now-duplicated old \"start screen\" bullet ---\nbody.paragraphs.load(\"items,text\");\nawait context.sync();\nlet startScreenIdxs = findAllIndices(body.paragraphs.items, \"start screen that allows\");\n// startScreenIdxs[0] is the bullet we just retargeted in Step 3 (keep it);\n// delete the later, now-stale duplicate.\nbody.paragraphs.items[startScreenIdxs[1]].delete();\nawait context.sync();\n\n// --- Step 6: \"Add a timer...\" -> \"Create a \\u201cfree play\\u201d mode...\" ---\nbody.paragraphs.load(\"items,text\");\nawait context.sync();\nidx = findIndex(body.paragraphs.items, \"Add a timer so the player\");\nbody.paragraphs.items[idx].insertOoxml(freePlayOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The \"Bonus:\" section contains a bulleted list of ideas. The edit:\n#  1. Turns the \"Create a 3+ frame explosion...\" bullet into \"Add the concept of lives\"\n#  2. Inserts a new bullet right after it with the old \"Alter the code...\" text\n#  3. Turns the (old) \"Add the concept of lives\" bullet into the \"start screen\" bullet\n#  4. Removes the now-duplicated old \"Alter the code...\" bullet\n#  5. Removes the now-duplicated old \"start screen\" bullet\n#  6. Turns the \"Add a timer...\" bullet into the new \"free play\" mode bullet\n#\n# Range.InsertXML(...) gives full control over the exact run/formatting\n# structure, matching the target XML precisely (mirrors Office.js'\n# insertOoxml). We always re-locate paragraphs by their (currently still\n# unique) text right before each such call so the script tolerates the\n# paragraph-count churn caused by the insert/delete steps.\n\n$d = $word.ActiveDocument\n\nfunction FindParaIndex($searchText) {\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        if ($d.Paragraphs.Item($i).Range.Text -like \"*$searchText*\") {\n            return $i\n        }\n    }\n    return -1\n}\n\nfunction FindAllParaIndices($searchText) {\n    $out = @()\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        if ($d.Paragraphs.Item($i).Range.Text -like \"*$searchText*\") {\n            $out += $i\n        }\n    }\n    return $out\n}\n\nfunction XmlWrap($innerParagraphXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $innerParagraphXml + '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n}\n\n$livesXml = XmlWrap('<w:p><w:pPr><w:ind w:firstLine=\"540\"/><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr></w:pPr><w:r><w:t xml:space=\"preserve\">\u25a1 </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\">Add </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\">the concept of </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>lives</w:t></w:r></w:p>')\n\n$shiesXml = XmlWrap('<w:p><w:pPr><w:ind w:firstLine=\"540\"/><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr></w:pPr><w:r><w:t>\u25a1</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>Alter the code so instead of following a line, the tracking bot &#8220;shies&#8221; away from a line</w:t></w:r></w:p>')\n\n$startScreenXml = XmlWrap('<w:p><w:pPr><w:ind w:firstLine=\"540\"/><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\">        </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t xml:space=\"preserve\">\u25a1 </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> Create</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> a start screen that allows the user to select which mode (follow or shy) to use</w:t></w:r></w:p>')\n\n$freePlayXml = XmlWrap('<w:p><w:pPr><w:ind w:firstLine=\"540\"/><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr></w:pPr><w:r><w:t xml:space=\"preserve\">\u25a1 </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>Create a &#8220;free play&#8221; mode where the player gets points for hitting random popup sprites</w:t></w:r></w:p>')\n\n# --- Step 1: \"Create a 3+ frame explosion...\" -> \"Add the concept of lives\" ---\n$idx = FindParaIndex(\"explosion animation\")\n$d.Paragraphs.Item($idx).Range.InsertXML($livesXml)\n\n# --- Step 2: insert an empty paragraph right after it ---\n$idx = FindParaIndex(\"the concept of lives\")\n$d.Paragraphs.Item($idx).Range.InsertParagraphAfter()\n\n# --- Step 3: give that new (now-following) paragraph the \"Alter the code...\" content ---\n$idx = FindParaIndex(\"the concept of lives\")\n$d.Paragraphs.Item($idx + 1).Range.InsertXML($shiesXml)\n\n# --- Step 4: turn the ORIGINAL \"Add the concept of lives\" bullet (further down) into the start-screen bullet ---\n$livesIdxs = FindAllParaIndices(\"the concept of lives\")\n# $livesIdxs[0] is the bullet created in Step 1 (keep it); the stale\n# original is the later duplicate.\n$d.Paragraphs.Item($livesIdxs[1]).Range.InsertXML($startScreenXml)\n\n# --- Step 5: delete the now-duplicated old \"Alter the code...\" bullet ---\n$shiesIdxs = FindAllParaIndices(\"shies\" + [char]0x201D + \" away from a line\")\n# $shiesIdxs[0] is the bullet created in Step 3 (keep it); delete the later,\n# now-stale duplicate.\n$d.Paragraphs.Item($shiesIdxs[1]).Range.Delete()\n\n# --- Step 6: delete the now-duplicated old \"start screen\" bullet ---\n$startScreenIdxs = FindAllParaIndices(\"start screen that allows\")\n# $startScreenIdxs[0] is the bullet just retargeted in Step 4 (keep it);\n# delete the later, now-stale duplicate.\n$d.Paragraphs.Item($startScreenIdxs[1]).Range.Delete()\n\n# --- Step 7: \"Add a timer...\" -> \"Create a \"free play\" mode...\" ---\n$idx = FindParaIndex(\"Add a timer so the player\")\n$d.Paragraphs.Item($idx).Range.InsertXML($freePlayXml)\n"}
